$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header row (row 1): tweak D1 text, add new F1 "Status" header with the
# same bold/centered/bordered look as the existing header cells.
# ---------------------------------------------------------------------------
$ws.Cells.Item(1, 4).Value = "Winning Numbers"

$ws.Cells.Item(1, 6).Value = "Status"
$ws.Cells.Item(1, 5).Copy()
$ws.Cells.Item(1, 6).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# The old rows 2-3 held numeric draw numbers / date-serial draw dates with a
# custom number format. The refreshed data is plain text end to end (draw
# numbers, dates and bonus balls are all stored as strings now). Write each
# value with a leading apostrophe so Excel stores it as literal text instead
# of re-parsing "1234"/"2025-05-01" back into a number or date, then strip
# the resulting quote-prefix/number formatting so the cells end up with
# plain, default (General) formatting - matching freshly authored text
# cells rather than ones carrying leftover numeric styling.
# ---------------------------------------------------------------------------
function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'" + $text
    $cell.ClearFormats()
}

# Row 2: Lotto
$ws.Cells.Item(2, 1).Value = "Lotto"
Set-TextCell 2 2 "1234"
Set-TextCell 2 3 "2025-05-01"
$ws.Cells.Item(2, 4).Value = "01, 05, 12, 26, 33, 45"
Set-TextCell 2 5 "22"
$ws.Cells.Item(2, 6).Value = "Missing"

# Row 3: Lotto Plus 1
$ws.Cells.Item(3, 1).Value = "Lotto Plus 1"
Set-TextCell 3 2 "567"
Set-TextCell 3 3 "2025-05-01"
$ws.Cells.Item(3, 4).Value = "07, 14, 21, 28, 35, 42"
Set-TextCell 3 5 "17"
$ws.Cells.Item(3, 6).Value = "Missing"

# Row 4: PowerBall (new row)
$ws.Cells.Item(4, 1).Value = "PowerBall"
Set-TextCell 4 2 "890"
Set-TextCell 4 3 "2025-05-02"
$ws.Cells.Item(4, 4).Value = "03, 11, 22, 33, 44"
Set-TextCell 4 5 "09"
$ws.Cells.Item(4, 6).Value = "Missing"

Write-Host "edit complete"
